$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 3
$ws.Range("G2").Value2 = 14.12404233333333
$ws.Range("H2").Value2 = 42.372127
$ws.Range("I2").Value2 = 0.8844735734357805
$ws.Range("J2").Value2 = 0.8844735734357805
$ws.Range("K2").Value2 = 3
$ws.Range("M2").Value2 = 0.8985896666666667
$ws.Range("N2").Value2 = 2.695769
$ws.Range("O2").Value2 = 0.150721683826239
$ws.Range("P2").Value2 = 0.150721683826239
$ws.Range("Q2").Value2 = 12.69171849229589
$ws.Range("R2").Value2 = 114.225466430663
$ws.Range("S2").Value2 = 0.1333093462880515
$ws.Range("T2").Value2 = 0.1333093462880515

$ws.Range("E3").Value2 = 3
$ws.Range("G3").Value2 = 14.12404233333333
$ws.Range("H3").Value2 = 42.372127
$ws.Range("I3").Value2 = 0.8844735734357805
$ws.Range("J3").Value2 = 0.8844735734357805
$ws.Range("K3").Value2 = 3
$ws.Range("M3").Value2 = 1.645054333333333
$ws.Range("N3").Value2 = 4.935163
$ws.Range("O3").Value2 = 0.2759272316422339
$ws.Range("P3").Value2 = 0.2759272316422339
$ws.Range("Q3").Value2 = 23.23481704463344
$ws.Range("R3").Value2 = 209.113353401701
$ws.Range("S3").Value2 = 0.2440503445788489
$ws.Range("T3").Value2 = 0.244050344578849

$ws.Range("E4").Value2 = 3
$ws.Range("G4").Value2 = 14.12404233333333
$ws.Range("H4").Value2 = 42.372127
$ws.Range("I4").Value2 = 0.8844735734357805
$ws.Range("J4").Value2 = 0.8844735734357805
$ws.Range("K4").Value2 = 3
$ws.Range("M4").Value2 = 3.418269666666667
$ws.Range("N4").Value2 = 10.254809
$ws.Range("O4").Value2 = 0.5733510845315271
$ws.Range("P4").Value2 = 0.5733510845315271
$ws.Range("Q4").Value2 = 48.27978547874923
$ws.Range("R4").Value2 = 434.518069308743
$ws.Range("S4").Value2 = 0.50711388256888
$ws.Range("T4").Value2 = 0.50711388256888

$ws.Range("E5").Value2 = 3
$ws.Range("G5").Value2 = 1.844826333333333
$ws.Range("H5").Value2 = 5.534479
$ws.Range("I5").Value2 = 0.1155264265642196
$ws.Range("J5").Value2 = 0.1155264265642196
$ws.Range("K5").Value2 = 3
$ws.Range("M5").Value2 = 0.8985896666666667
$ws.Range("N5").Value2 = 2.695769
$ws.Range("O5").Value2 = 0.150721683826239
$ws.Range("P5").Value2 = 0.150721683826239
$ws.Range("Q5").Value2 = 1.657741879927889
$ws.Range("R5").Value2 = 14.919676919351
$ws.Range("S5").Value2 = 0.01741233753818752
$ws.Range("T5").Value2 = 0.01741233753818753

$ws.Range("E6").Value2 = 3
$ws.Range("G6").Value2 = 1.844826333333333
$ws.Range("H6").Value2 = 5.534479
$ws.Range("I6").Value2 = 0.1155264265642196
$ws.Range("J6").Value2 = 0.1155264265642196
$ws.Range("K6").Value2 = 3
$ws.Range("M6").Value2 = 1.645054333333333
$ws.Range("N6").Value2 = 4.935163
$ws.Range("O6").Value2 = 0.2759272316422339
$ws.Range("P6").Value2 = 0.2759272316422339
$ws.Range("Q6").Value2 = 3.034839553897445
$ws.Range("R6").Value2 = 27.313555985077
$ws.Range("S6").Value2 = 0.03187688706338493
$ws.Range("T6").Value2 = 0.03187688706338494

$ws.Range("E7").Value2 = 3
$ws.Range("G7").Value2 = 1.844826333333333
$ws.Range("H7").Value2 = 5.534479
$ws.Range("I7").Value2 = 0.1155264265642196
$ws.Range("J7").Value2 = 0.1155264265642196
$ws.Range("K7").Value2 = 3
$ws.Range("M7").Value2 = 3.418269666666667
$ws.Range("N7").Value2 = 10.254809
$ws.Range("O7").Value2 = 0.5733510845315271
$ws.Range("P7").Value2 = 0.5733510845315271
$ws.Range("Q7").Value2 = 6.306113895501223
$ws.Range("R7").Value2 = 56.755025059511
$ws.Range("S7").Value2 = 0.0662372019626471
$ws.Range("T7").Value2 = 0.0662372019626471
